# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-board figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) for a handful of leve rows across the per-job
# sheets, same as the scheduled price-refresh job would. Values only -
# no structural changes.

$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values   # column letter -> new value ($null clears the cell)
    )

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $cell = $ws.Range("$col$Row")
        if ($null -eq $Values[$col]) {
            $cell.ClearContents()
        } else {
            $cell.Value = $Values[$col]
        }
    }
}

# ALC
Set-LeveRow "ALC" 130 @{
    H = 52720; J = 52720; L = 52720; N = -62760
}
Set-LeveRow "ALC" 132 @{
    H = 1197.341; I = 876.2368; J = 3231; K = 2628.7104; L = 9693
    M = -98.71039999999994; N = -14753
}
Set-LeveRow "ALC" 137 @{
    H = 1434.2963; I = 1109.975; K = 3329.925; M = -779.9249999999997
}
Set-LeveRow "ALC" 139 @{
    H = 44853.332; J = 44853.332; L = 44853.332; N = -55133.332
}

# ARM
Set-LeveRow "ARM" 6 @{
    H = 0; I = 0; J = 0; K = 0; L = 0; M = $null; N = $null
}
Set-LeveRow "ARM" 32 @{
    H = 1235.09; I = 1147; J = 1881.0834; K = 1147; L = 1881.0834
    M = -860; N = -2455.0834
}
Set-LeveRow "ARM" 61 @{
    H = 3283.9583; I = 3649.8975; J = 1698.2222; K = 3649.8975; L = 1698.2222
    M = -3437.8975; N = -2122.2222
}
Set-LeveRow "ARM" 74 @{
    H = 1177.2709; I = 1012.5476; J = 2330.3333; K = 1012.5476; L = 2330.3333
    M = -138.5476; N = -4078.3333
}
Set-LeveRow "ARM" 77 @{
    H = 1177.2709; I = 1012.5476; J = 2330.3333; K = 5062.738; L = 11651.6665
    M = -694.7380000000003; N = -20387.6665
}
Set-LeveRow "ARM" 122 @{
    H = 1710204.9; I = 2137488.2; J = 1071; K = 6412464.600000001; L = 3213
    M = -6410014.600000001; N = -8113
}
Set-LeveRow "ARM" 132 @{
    H = 2259.0469; I = 1468.8654; J = 5683.1665; K = 4406.5962; L = 17049.4995
    M = -1876.5962; N = -22109.4995
}
Set-LeveRow "ARM" 136 @{
    H = 3283.9583; I = 3649.8975; J = 1698.2222; K = 10949.6925; L = 5094.6666
    M = -8399.692500000001; N = -10194.6666
}

# CRP
Set-LeveRow "CRP" 31 @{
    H = 224011.34; I = 1622.5333; J = 569097.4399999999; K = 1622.5333
    L = 569097.4399999999; M = -1327.5333; N = -569687.4399999999
}
Set-LeveRow "CRP" 34 @{
    H = 224011.34; I = 1622.5333; J = 569097.4399999999; K = 1622.5333
    L = 569097.4399999999; M = -1420.5333; N = -569501.4399999999
}
Set-LeveRow "CRP" 134 @{
    H = 1391.9036; I = 1450.305; K = 4350.915; M = -1815.915
}

# CUL
Set-LeveRow "CUL" 55 @{
    H = 2370.28; J = 2370.28; L = 7110.84; N = -7464.84
}

# GSM
Set-LeveRow "GSM" 3 @{
    H = 3000; I = 3000; J = 0; K = 3000; L = 0; M = -2884; N = $null
}
Set-LeveRow "GSM" 102 @{
    H = 514801.25; I = 652813.1; K = 652813.1; M = -651191.1
}
Set-LeveRow "GSM" 122 @{
    H = 25972548; I = 42593820; J = 1813.8125; K = 127781460; L = 5441.4375
    M = -127779010; N = -10341.4375
}
Set-LeveRow "GSM" 126 @{
    H = 4038.7646; I = 5390.3076; J = 2633.16; K = 16170.9228; L = 7899.48
    M = -13700.9228; N = -12839.48
}
Set-LeveRow "GSM" 132 @{
    H = 1495.6028; I = 1085.1404; J = 2957.875; K = 3255.4212; L = 8873.625
    M = -725.4211999999998; N = -13933.625
}

# LTW
Set-LeveRow "LTW" 22 @{
    H = 5208982; I = 10417025; J = 939.25; K = 10417025; L = 939.25
    M = -10416730; N = -1529.25
}
Set-LeveRow "LTW" 27 @{
    H = 5208982; I = 10417025; J = 939.25; K = 10417025; L = 939.25
    M = -10416918; N = -1153.25
}
Set-LeveRow "LTW" 132 @{
    H = 12147812; I = 15719910; J = 2679.8; K = 47159730; L = 8039.400000000001
    M = -47157200; N = -13099.4
}
Set-LeveRow "LTW" 136 @{
    H = 6632.6665; I = 4282.4287; K = 12847.2861; M = -10297.2861
}

# WVR
Set-LeveRow "WVR" 107 @{
    H = 80000620; I = 142857940; J = 6667100.5; K = 428573820; L = 20001301.5
    M = -428571900; N = -20005141.5
}
Set-LeveRow "WVR" 132 @{
    H = 17624.217; I = 19928.27; J = 2647.875; K = 59784.81; L = 7943.625
    M = -57254.81; N = -13003.625
}
Set-LeveRow "WVR" 136 @{
    H = 8774624; I = 2923.9722; J = 23811824; K = 8771.9166; L = 71435472
    M = -6221.9166; N = -71440572
}
